$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041718129476462
$ws.Range("D2").Value = 1.042656494379526
$ws.Range("E2").Value = 1.039923658474874
$ws.Range("F2").Value = 1.049810263742105
$ws.Range("I2").Value = 1.040884733319624
$ws.Range("J2").Value = 1.04679782858849
$ws.Range("K2").Value = 1.045432565200145
$ws.Range("L2").Value = 1.042707466958246
$ws.Range("M2").Value = 1.052566286466878
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042869221613038
$ws.Range("D3").Value = 1.043516771764116
$ws.Range("E3").Value = 1.040909820301119
$ws.Range("F3").Value = 1.051165926162506
$ws.Range("I3").Value = 1.041227476843246
$ws.Range("J3").Value = 1.047594060588326
$ws.Range("K3").Value = 1.046103589213633
$ws.Range("L3").Value = 1.043503483737747
$ws.Range("M3").Value = 1.05373286805462
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04361378452914
$ws.Range("D4").Value = 1.044073158814578
$ws.Range("E4").Value = 1.04154798876888
$ws.Range("F4").Value = 1.05204321100008
$ws.Range("I4").Value = 1.041447928910683
$ws.Range("J4").Value = 1.048108469608141
$ws.Range("K4").Value = 1.046536877311261
$ws.Range("L4").Value = 1.044018003664952
$ws.Range("M4").Value = 1.054487269802906
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043926735579288
$ws.Range("D5").Value = 1.044306999761361
$ws.Range("E5").Value = 1.041816288834822
$ws.Range("F5").Value = 1.052412043395755
$ws.Range("I5").Value = 1.041540289945504
$ws.Range("J5").Value = 1.048324534994144
$ws.Range("K5").Value = 1.046718814549536
$ws.Range("L5").Value = 1.044234175411174
$ws.Range("M5").Value = 1.054804313639439
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04397927777807
$ws.Range("D6").Value = 1.044346258935551
$ws.Range("E6").Value = 1.041861338439575
$ws.Range("F6").Value = 1.052473973333474
$ws.Range("I6").Value = 1.041555779191745
$ws.Range("J6").Value = 1.048360802079759
$ws.Range("K6").Value = 1.046749349900623
$ws.Range("L6").Value = 1.044270463853463
$ws.Range("M6").Value = 1.054857540502967
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043617966442717
$ws.Range("D7").Value = 1.044076283660331
$ws.Range("E7").Value = 1.041551573751472
$ws.Range("F7").Value = 1.052048139267552
$ws.Range("I7").Value = 1.041449164288486
$ws.Range("J7").Value = 1.048111357438399
$ws.Range("K7").Value = 1.04653930921775
$ws.Range("L7").Value = 1.04402089268188
$ws.Range("M7").Value = 1.054491506574561
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042107203459481
$ws.Range("D8").Value = 1.042947285318417
$ws.Range("E8").Value = 1.040256925331787
$ws.Range("F8").Value = 1.05026840105257
$ws.Range("I8").Value = 1.041000839884692
$ws.Range("J8").Value = 1.047067086400659
$ws.Range("K8").Value = 1.04565952955525
$ws.Range("L8").Value = 1.04297660014253
$ws.Range("M8").Value = 1.052960633593913
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039442901178528
$ws.Range("D9").Value = 1.040955757544342
$ws.Range("E9").Value = 1.037975986035269
$ws.Range("F9").Value = 1.047132777763804
$ws.Range("I9").Value = 1.040200661456431
$ws.Range("J9").Value = 1.045220732527836
$ws.Range("K9").Value = 1.044102253303261
$ws.Range("L9").Value = 1.041132131337511
$ws.Range("M9").Value = 1.050259461606173
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037665154636979
$ws.Range("D10").Value = 1.039626636074162
$ws.Range("E10").Value = 1.036455580531095
$ws.Range("F10").Value = 1.045042538237022
$ws.Range("I10").Value = 1.039660341656542
$ws.Range("J10").Value = 1.043985594087993
$ws.Range("K10").Value = 1.043059326951387
$ws.Range("L10").Value = 1.039899552159258
$ws.Range("M10").Value = 1.048456136972074
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036894980635277
$ws.Range("D11").Value = 1.03905076312287
$ws.Range("E11").Value = 1.035797269413578
$ws.Range("F11").Value = 1.04413744648245
$ws.Range("I11").Value = 1.039424741811675
$ws.Range("J11").Value = 1.043449747347724
$ws.Range("K11").Value = 1.042606593715038
$ws.Range("L11").Value = 1.039365124687688
$ws.Range("M11").Value = 1.047674643014217
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036608841783192
$ws.Range("D12").Value = 1.038836804204969
$ws.Range("E12").Value = 1.035552747481125
$ws.Range("F12").Value = 1.043801251217015
$ws.Range("I12").Value = 1.039336982850148
$ws.Range("J12").Value = 1.043250554964336
$ws.Range("K12").Value = 1.042438256338508
$ws.Range("L12").Value = 1.039166506407651
$ws.Range("M12").Value = 1.047384262052373
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0366702223754
$ws.Range("D13").Value = 1.038882701578389
$ws.Range("E13").Value = 1.035605198095895
$ws.Range("F13").Value = 1.04387336647093
$ws.Range("I13").Value = 1.039355818618587
$ws.Range("J13").Value = 1.043293289443071
$ws.Range("K13").Value = 1.042474373079928
$ws.Range("L13").Value = 1.039209115612889
$ws.Range("M13").Value = 1.047446554293138
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036871329584358
$ws.Range("D14").Value = 1.039033078319436
$ws.Range("E14").Value = 1.035777057076277
$ws.Range("F14").Value = 1.044109656575084
$ws.Range("I14").Value = 1.039417492664553
$ws.Range("J14").Value = 1.043433285213326
$ws.Range("K14").Value = 1.042592682411087
$ws.Range("L14").Value = 1.039348709050548
$ws.Range("M14").Value = 1.047650642075333
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036995230102035
$ws.Range("D15").Value = 1.03912572320502
$ws.Range("E15").Value = 1.035882945581481
$ws.Range("F15").Value = 1.044255242094187
$ws.Range("I15").Value = 1.039455459372294
$ws.Range("J15").Value = 1.043519520645198
$ws.Range("K15").Value = 1.04266555386244
$ws.Range("L15").Value = 1.039434702815437
$ws.Range("M15").Value = 1.04777637405083
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037716259980787
$ws.Range("D16").Value = 1.039664847327306
$ws.Range("E16").Value = 1.036499271108785
$ws.Range("F16").Value = 1.045102605829156
$ws.Range("I16").Value = 1.039675943082119
$ws.Range("J16").Value = 1.044021134797382
$ws.Range("K16").Value = 1.043089349282877
$ws.Range("L16").Value = 1.03993500524666
$ws.Range("M16").Value = 1.048507988378315
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038168434906374
$ws.Range("D17").Value = 1.040002930157531
$ws.Range("E17").Value = 1.036885884322609
$ws.Range("F17").Value = 1.045634131284458
$ws.Range("I17").Value = 1.039813807758907
$ws.Range("J17").Value = 1.044335509160019
$ws.Range("K17").Value = 1.043354879286512
$ws.Range("L17").Value = 1.040248640220612
$ws.Range("M17").Value = 1.048966736741852
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038432142631035
$ws.Range("D18").Value = 1.040200093993694
$ws.Range("E18").Value = 1.037111392669022
$ws.Range("F18").Value = 1.04594416089859
$ws.Range("I18").Value = 1.039894063910019
$ws.Range("J18").Value = 1.044518779724436
$ws.Range("K18").Value = 1.043509648587751
$ws.Range("L18").Value = 1.040431509515493
$ws.Range("M18").Value = 1.049234255096799
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038522053725857
$ws.Range("D19").Value = 1.040267315996692
$ws.Range("E19").Value = 1.037188285796155
$ws.Range("F19").Value = 1.046049873051462
$ws.Range("I19").Value = 1.039921402406051
$ws.Range("J19").Value = 1.044581253550586
$ws.Range("K19").Value = 1.043562402314555
$ws.Range("L19").Value = 1.040493851602539
$ws.Range("M19").Value = 1.049325461544058
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038119924788691
$ws.Range("D20").Value = 1.039966660608169
$ws.Range("E20").Value = 1.036844404039401
$ws.Range("F20").Value = 1.045577103676293
$ws.Range("I20").Value = 1.039799032529107
$ws.Range("J20").Value = 1.04430178998217
$ws.Range("K20").Value = 1.043326401799737
$ws.Range("L20").Value = 1.040214997267454
$ws.Range("M20").Value = 1.048917523796446
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036812110250713
$ws.Range("D21").Value = 1.038988797616025
$ws.Range("E21").Value = 1.035726448784579
$ws.Range("F21").Value = 1.04404007516227
$ws.Range("I21").Value = 1.039399338010071
$ws.Range("J21").Value = 1.043392064238398
$ws.Range("K21").Value = 1.042557848020035
$ws.Range("L21").Value = 1.039307605255725
$ws.Range("M21").Value = 1.047590546074047
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035989475939894
$ws.Range("D22").Value = 1.038373663153638
$ws.Range("E22").Value = 1.035023569300584
$ws.Range("F22").Value = 1.043073660218563
$ws.Range("I22").Value = 1.039146606623776
$ws.Range("J22").Value = 1.042819185886846
$ws.Range("K22").Value = 1.042073631602187
$ws.Range("L22").Value = 1.03873646579238
$ws.Range("M22").Value = 1.046755646554071
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036425604735592
$ws.Range("D23").Value = 1.03869978752324
$ws.Range("E23").Value = 1.035396177154049
$ws.Range("F23").Value = 1.043585978438557
$ws.Range("I23").Value = 1.03928071984131
$ws.Range("J23").Value = 1.043122965049098
$ws.Range("K23").Value = 1.042330418624534
$ws.Range("L23").Value = 1.039039297332619
$ws.Range("M23").Value = 1.047198298111684
$ws.Range("B24").Value = 1.019999999999999
$ws.Range("C24").Value = 1.038141844540425
$ws.Range("D24").Value = 1.03998304936158
$ws.Range("E24").Value = 1.036863147183038
$ws.Range("F24").Value = 1.045602871997485
$ws.Range("I24").Value = 1.039805709307401
$ws.Range("J24").Value = 1.044317026531754
$ws.Range("K24").Value = 1.043339269889353
$ws.Range("L24").Value = 1.040230199281715
$ws.Range("M24").Value = 1.048939761198386
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040131951944615
$ws.Range("D25").Value = 1.041470865989525
$ws.Range("E25").Value = 1.038565621229035
$ws.Range("F25").Value = 1.047943368628281
$ws.Range("I25").Value = 1.040408735193858
$ws.Range("J25").Value = 1.045698801643551
$ws.Range("K25").Value = 1.044505679357241
$ws.Range("L25").Value = 1.041609484241801
$ws.Range("M25").Value = 1.050958217493167

Write-Output "Updated 240 cells in rows 2-25 (B,C,D,E,F,I,J,K,L,M) with new power flow results for 380 kV case."
